# Add most recent smelt data (Week 23 / November 3-7, 2025)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 - per-stratum weekly abundance data
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 - Week number -> Dates lookup

# ---------------------------------------------------------------------
# Sheet1: append week 23 rows (151-157), one per stratum plus the
# "All Strata" totals row, mirroring the layout of the preceding weeks.
# ---------------------------------------------------------------------

$week23 = @(
    @{ Row = 151; Stratum = "Suisun Bay";          Sites = 3;  Tows = 12;  Caught = 0;  Marked = 0; Code9 = 0; Volume = 50456;  Index = "0*";  Lower = "NA";  Upper = "NA" },
    @{ Row = 152; Stratum = "Suisun Marsh";         Sites = 6;  Tows = 20;  Caught = 19; Marked = 0; Code9 = 0; Volume = 89408;  Index = 7125;  Lower = 1518;  Upper = 21253 },
    @{ Row = 153; Stratum = "Lower Sacramento";      Sites = 6;  Tows = 24;  Caught = 0;  Marked = 0; Code9 = 0; Volume = 93981;  Index = "0*";  Lower = "NA";  Upper = "NA" },
    @{ Row = 154; Stratum = "Cache Slough LI";        Sites = 3;  Tows = 12;  Caught = 0;  Marked = 0; Code9 = 0; Volume = 40273;  Index = "0*";  Lower = "NA";  Upper = "NA" },
    @{ Row = 155; Stratum = "Sac DW Ship Channel";    Sites = 6;  Tows = 24;  Caught = 0;  Marked = 0; Code9 = 0; Volume = 87009;  Index = "0*";  Lower = "NA";  Upper = "NA" },
    @{ Row = 156; Stratum = "Lower San Joaquin";      Sites = 6;  Tows = 24;  Caught = 0;  Marked = 0; Code9 = 0; Volume = 91437;  Index = "0*";  Lower = "NA";  Upper = "NA" },
    @{ Row = 157; Stratum = "All Strata";             Sites = 30; Tows = 116; Caught = 19; Marked = 0; Code9 = 0; Volume = 452563; Index = 7125;  Lower = 1518;  Upper = 21253 }
)

foreach ($r in $week23) {
    $row = $r.Row
    $ws1.Cells.Item($row, 1).Value = 23
    $ws1.Cells.Item($row, 2).Value = $r.Stratum
    $ws1.Cells.Item($row, 3).Value = $r.Sites
    $ws1.Cells.Item($row, 4).Value = $r.Tows
    $ws1.Cells.Item($row, 5).Value = $r.Caught
    $ws1.Cells.Item($row, 6).Value = $r.Marked
    $ws1.Cells.Item($row, 7).Value = $r.Code9
    $ws1.Cells.Item($row, 8).Value = $r.Volume
    $ws1.Cells.Item($row, 8).NumberFormat = "#,##0"
    $ws1.Cells.Item($row, 9).Value = $r.Index
    $ws1.Cells.Item($row, 10).Value = $r.Lower
    $ws1.Cells.Item($row, 11).Value = $r.Upper

    if ($r.Index -is [int]) {
        $ws1.Cells.Item($row, 9).NumberFormat = "#,##0"
    }
    if ($r.Lower -is [int]) {
        $ws1.Cells.Item($row, 10).NumberFormat = "#,##0"
    }
    if ($r.Upper -is [int]) {
        $ws1.Cells.Item($row, 11).NumberFormat = "#,##0"
    }
}

# ---------------------------------------------------------------------
# Sheet2: append the new Week 23 / date lookup row.
# ---------------------------------------------------------------------

$ws2.Cells.Item(24, 1).Value = 23
$ws2.Cells.Item(24, 2).Value = "November 3–7, 2025"

# ---------------------------------------------------------------------
# View state: the edit was made while looking at Sheet1 (F140), leaving
# Sheet2 unselected.
# ---------------------------------------------------------------------

$ws2.Range("F20:F21").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("F140").Select() | Out-Null
